$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values look numeric to Excel auto-detection (e.g. "1.010", "313.30")
# and would be silently coerced to plain numbers (losing formatting / becoming t="n")
# if assigned directly via .Value. To preserve them as literal text (matching the
# original inlineStr cells) we write each one as a text-returning formula, then
# convert the whole range to static values via Copy + PasteSpecial(xlPasteValues).

$ws.Range("D2").Formula = '="28.129.55"'
$ws.Range("D3").Formula = '="1.850.90"'
$ws.Range("D4").Formula = '="1.010"'
$ws.Range("D5").Formula = '="313.30"'
$ws.Range("D6").Formula = '="1.009"'
$ws.Range("D7").Formula = '="0.5079"'
$ws.Range("D8").Formula = '="0.3899"'
$ws.Range("D9").Formula = '="0.08232"'
$ws.Range("D10").Formula = '="1.106"'
$ws.Range("D11").Formula = '="6.178"'
$ws.Range("D12").Formula = '="1.856.86"'
$ws.Range("D13").Formula = '="20.14"'
$ws.Range("D14").Formula = '="7.156"'
$ws.Range("D15").Formula = '="1.008"'
$ws.Range("D16").Formula = '="0.00001095"'
$ws.Range("D17").Formula = '="90.87"'
$ws.Range("D18").Formula = '="0.06674"'
$ws.Range("D19").Formula = '="1.009"'
$ws.Range("D20").Formula = '="17.47"'
$ws.Range("D21").Formula = '="5.904"'
$ws.Range("D22").Formula = '="28.165.16"'
$ws.Range("D23").Formula = '="11.00"'
$ws.Range("D24").Formula = '="2.239"'
$ws.Range("D25").Formula = '="2.068.85"'
$ws.Range("D26").Formula = '="159.73"'
$ws.Range("D27").Formula = '="20.57"'
$ws.Range("D28").Formula = '="2.368"'
$ws.Range("D29").Formula = '="125.80"'
$ws.Range("D30").Formula = '="0.1034"'
$ws.Range("D31").Formula = '="1.018"'
$ws.Range("D32").Formula = '="5.759"'
$ws.Range("D33").Formula = '="3.619"'
$ws.Range("D34").Formula = '="0.02407"'
$ws.Range("D35").Formula = '="0.06405"'
$ws.Range("D36").Formula = '="9.046"'
$ws.Range("D37").Formula = '="0.2159"'
$ws.Range("D38").Formula = '="1.245"'
$ws.Range("D39").Formula = '="1.172"'
$ws.Range("D40").Formula = '="0.6364"'
$ws.Range("D41").Formula = '="4.907"'
$ws.Range("D42").Formula = '="11.02"'
$ws.Range("D43").Formula = '="0.5948"'
$ws.Range("D44").Formula = '="12.79"'
$ws.Range("D45").Formula = '="3.676"'
$ws.Range("D46").Formula = '="1.275"'
$ws.Range("D47").Formula = '="1.957"'
$ws.Range("D48").Formula = '="1.192"'
$ws.Range("D49").Formula = '="120.10"'
$ws.Range("D50").Formula = '="0.06809"'
$ws.Range("D51").Formula = '="75.68"'

$priceRange = $ws.Range("D2:D51")
$priceRange.Copy()
$priceRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Column E (Volume/1h change) values are padded percentage strings (e.g. "  -0.53%  ")
# which Excel does not auto-convert to numbers, so a direct text assignment is safe.

$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  -6.00%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("E32").Value = "  -2.32%  "
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("E36").Value = "  -7.79%  "
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("E44").Value = "  -2.41%  "
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("E51").Value = "  -4.09%  "
